# Add a new "2020" column (J) to the Sustainable Development Goal 12.4.2
# indicator table, matching the extra year of data added upstream, then
# drop the two now-unused trailing blank rows (27 and 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J3: year header "2020" (same look as the other year header cells) ---
$ws.Range("I3").Copy($ws.Range("J3"))
$ws.Range("J3").Value = 2020

# --- J4: "Hazardous waste generation (1-3 hazard classes), thousand tons" value ---
$ws.Range("I4").Copy($ws.Range("J4"))
$ws.Range("J4").Value = 11545.7
$ws.Range("J4").NumberFormat = "0.0"

# --- J5: "Resident population, thousand people" value (reported as text) ---
$ws.Range("I5").Copy($ws.Range("J5"))
$ws.Range("J5").Value = "1 754,6"
$ws.Range("J5").NumberFormat = "0.0"
$ws.Range("J5").HorizontalAlignment = -4152

# --- J6: "Generation of hazardous waste per person, kilogram / person" value ---
$ws.Range("I6").Copy($ws.Range("J6"))
$ws.Range("J6").Value = 6636.8

# Remove the two trailing empty rows that are no longer part of the table.
$ws.Rows("27:28").Delete()

# Leave the same cell selected as in the published workbook.
$ws.Range("G22").Select()
